$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the ID column keeps storing text values (not auto-converted to numbers)
$ws.Range("E2:E3").NumberFormat = "@"

# Row 2 updates
$ws.Range("B2").Value = 2274
$ws.Range("C2").Value = "2022-11-14 15:00"
$ws.Range("E2").Value = "4542"

# Row 3 updates
$ws.Range("B3").Value = 2865
$ws.Range("C3").Value = "2022-11-14 15:00"
$ws.Range("E3").Value = "1964"
